$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new row at position 7; existing rows 7-10 shift down to 8-11.
$ws.Rows.Item(7).Insert()

# Select the newly inserted row (mirrors the author selecting/highlighting row 7).
$ws.Rows.Item(7).Select()

# Fill in the new organization's data in row 7.
$ws.Range("A7").Value = "Poudre Wilderness Volunteers"
$ws.Range("B7").Value = "Nonprofit"
$ws.Range("C7").Value = "Caring for Northern Colorado Wilderness."
$ws.Range("D7").Value = "Hiking, riding, educating."
$ws.Range("E7").Value = "https://www.pwv.org/"
$ws.Range("E7").Style = "Hyperlink"
$ws.Range("F7").Value = "Yes"

$ws.Range("G7:J7").NumberFormat = "0.000000"
$ws.Range("G7").Value = -105.07315
$ws.Range("H7").Value = 40.53125

$ws.Range("L7").Value = "PO Box 271921, Fort Collins, CO  80527"
$ws.Range("L7").Font.Name = "Trebuchet MS"
$ws.Range("L7").Font.Size = 10
$ws.Range("L7").Font.Color = 7368816

# Match the row height used for the new entry.
$ws.Rows.Item(7).RowHeight = 15
